function Set-CellAndLink($ws, $addr, $value) {
    $ws.Range($addr).Value = $value
    $target = $ws.Range($addr).Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            $hl.TextToDisplay = $value
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: the two tracked source files (b1beacd1... and f57829ef...)
# swap rows -- f57829ef (already handed back) now appears first (row 2), and
# b1beacd1 moves to row 3 with its status updated from "Handed back: in sync
# with en-US" to "Ready for handoff".
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Set-CellAndLink $wsOverview "A2" "f57829ef-8225-4cd2-890a-00cf54974452.md"

Set-CellAndLink $wsOverview "A3" "b1beacd1-6e4f-45bc-8352-436126d411ef.md"
Set-CellAndLink $wsOverview "B3" "Ready for handoff"
Set-CellAndLink $wsOverview "C3" "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-CellAndLink $wsZhCn "A2" "f57829ef-8225-4cd2-890a-00cf54974452.md"
Set-CellAndLink $wsZhCn "C2" "f57829ef-8225-4cd2-890a-00cf54974452.5029642f4f9b71fb3c5ba5befa9b2878be2109e9.zh-cn.xlf"
Set-CellAndLink $wsZhCn "E2" "f57829ef-8225-4cd2-890a-00cf54974452.md"
Set-CellAndLink $wsZhCn "F2" "f57829ef-8225-4cd2-890a-00cf54974452.5029642f4f9b71fb3c5ba5befa9b2878be2109e9.zh-cn.xlf"

Set-CellAndLink $wsZhCn "A3" "b1beacd1-6e4f-45bc-8352-436126d411ef.md"
Set-CellAndLink $wsZhCn "B3" "Ready for handoff"
Set-CellAndLink $wsZhCn "C3" "b1beacd1-6e4f-45bc-8352-436126d411ef.786bfc1b6fdee835f5a90e03138bb0dbc4f3f712.zh-cn.xlf"
$wsZhCn.Range("D3").Value = "2016-03-08 06:36:09"
Set-CellAndLink $wsZhCn "E3" "b1beacd1-6e4f-45bc-8352-436126d411ef.md"
Set-CellAndLink $wsZhCn "F3" "b1beacd1-6e4f-45bc-8352-436126d411ef.786bfc1b6fdee835f5a90e03138bb0dbc4f3f712.zh-cn.xlf"

# ---------------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-CellAndLink $wsDeDe "A2" "f57829ef-8225-4cd2-890a-00cf54974452.md"
Set-CellAndLink $wsDeDe "C2" "f57829ef-8225-4cd2-890a-00cf54974452.5029642f4f9b71fb3c5ba5befa9b2878be2109e9.de-de.xlf"
Set-CellAndLink $wsDeDe "E2" "f57829ef-8225-4cd2-890a-00cf54974452.md"
Set-CellAndLink $wsDeDe "F2" "f57829ef-8225-4cd2-890a-00cf54974452.5029642f4f9b71fb3c5ba5befa9b2878be2109e9.de-de.xlf"

Set-CellAndLink $wsDeDe "A3" "b1beacd1-6e4f-45bc-8352-436126d411ef.md"
Set-CellAndLink $wsDeDe "B3" "Ready for handoff"
Set-CellAndLink $wsDeDe "C3" "b1beacd1-6e4f-45bc-8352-436126d411ef.786bfc1b6fdee835f5a90e03138bb0dbc4f3f712.de-de.xlf"
$wsDeDe.Range("D3").Value = "2016-03-08 06:36:19"
Set-CellAndLink $wsDeDe "E3" "b1beacd1-6e4f-45bc-8352-436126d411ef.md"
Set-CellAndLink $wsDeDe "F3" "b1beacd1-6e4f-45bc-8352-436126d411ef.786bfc1b6fdee835f5a90e03138bb0dbc4f3f712.de-de.xlf"
